$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in M963 = 0 (previously-missing cell in last existing row)
$ws.Range("M963").Value = 0

# Append 30 new rows (964-993) of COVID data, one per day from 2022/11/17 to 2022/12/16
$dates = @(
  "2022/11/17", "2022/11/18", "2022/11/19", "2022/11/20", "2022/11/21", "2022/11/22", "2022/11/23", "2022/11/24", "2022/11/25", "2022/11/26", "2022/11/27", "2022/11/28", "2022/11/29", "2022/11/30", "2022/12/01", "2022/12/02", "2022/12/03", "2022/12/04", "2022/12/05", "2022/12/06", "2022/12/07", "2022/12/08", "2022/12/09", "2022/12/10", "2022/12/11", "2022/12/12", "2022/12/13", "2022/12/14", "2022/12/15", "2022/12/16"
)

# Row 964
$ws.Range("A964").Value = "'" + $dates[0]
$ws.Range("A964").Style = "Normal"
$ws.Range("B964").Value = 963
$ws.Range("C964").Value = 39486
$ws.Range("D964").Value = 691
$ws.Range("F964").Value = 0.01749987337
$ws.Range("G964").Value = 38709
$ws.Range("H964").Value = 27040
$ws.Range("I964").Value = 66526
$ws.Range("J964").Value = 38
$ws.Range("K964").Value = 0
$ws.Range("L964").Value = 38
$ws.Range("M964").Value = 0
$ws.Range("O964").Value = 0
$ws.Range("P964").Value = 1
$ws.Range("Q964").Value = 1
$ws.Range("T964").Value = 138

# Row 965
$ws.Range("A965").Value = "'" + $dates[1]
$ws.Range("A965").Style = "Normal"
$ws.Range("B965").Value = 964
$ws.Range("C965").Value = 39486
$ws.Range("D965").Value = 691
$ws.Range("F965").Value = 0.01749987337
$ws.Range("G965").Value = 38709
$ws.Range("H965").Value = 27040
$ws.Range("I965").Value = 66526
$ws.Range("J965").Value = 38
$ws.Range("K965").Value = 0
$ws.Range("L965").Value = 38
$ws.Range("M965").Value = 0
$ws.Range("O965").Value = 0
$ws.Range("P965").Value = 1
$ws.Range("Q965").Value = 1
$ws.Range("T965").Value = 138

# Row 966
$ws.Range("A966").Value = "'" + $dates[2]
$ws.Range("A966").Style = "Normal"
$ws.Range("B966").Value = 965
$ws.Range("C966").Value = 39486
$ws.Range("D966").Value = 691
$ws.Range("F966").Value = 0.01749987337
$ws.Range("G966").Value = 38709
$ws.Range("H966").Value = 27040
$ws.Range("I966").Value = 66526
$ws.Range("J966").Value = 38
$ws.Range("K966").Value = 0
$ws.Range("L966").Value = 38
$ws.Range("M966").Value = 0
$ws.Range("O966").Value = 0
$ws.Range("P966").Value = 1
$ws.Range("Q966").Value = 1
$ws.Range("T966").Value = 138

# Row 967
$ws.Range("A967").Value = "'" + $dates[3]
$ws.Range("A967").Style = "Normal"
$ws.Range("B967").Value = 966
$ws.Range("C967").Value = 39486
$ws.Range("D967").Value = 691
$ws.Range("F967").Value = 0.01749987337
$ws.Range("G967").Value = 38709
$ws.Range("H967").Value = 27040
$ws.Range("I967").Value = 66526
$ws.Range("J967").Value = 38
$ws.Range("K967").Value = 0
$ws.Range("L967").Value = 38
$ws.Range("M967").Value = 0
$ws.Range("O967").Value = 0
$ws.Range("P967").Value = 1
$ws.Range("Q967").Value = 1
$ws.Range("T967").Value = 139

# Row 968
$ws.Range("A968").Value = "'" + $dates[4]
$ws.Range("A968").Style = "Normal"
$ws.Range("B968").Value = 967
$ws.Range("C968").Value = 39486
$ws.Range("D968").Value = 691
$ws.Range("F968").Value = 0.01749987337
$ws.Range("G968").Value = 38709
$ws.Range("H968").Value = 27040
$ws.Range("I968").Value = 66526
$ws.Range("J968").Value = 38
$ws.Range("K968").Value = 0
$ws.Range("L968").Value = 38
$ws.Range("M968").Value = 0
$ws.Range("O968").Value = 0
$ws.Range("P968").Value = 1
$ws.Range("Q968").Value = 1
$ws.Range("T968").Value = 139

# Row 969
$ws.Range("A969").Value = "'" + $dates[5]
$ws.Range("A969").Style = "Normal"
$ws.Range("B969").Value = 968
$ws.Range("C969").Value = 39486
$ws.Range("D969").Value = 691
$ws.Range("F969").Value = 0.01749987337
$ws.Range("G969").Value = 38709
$ws.Range("H969").Value = 27040
$ws.Range("I969").Value = 66526
$ws.Range("J969").Value = 38
$ws.Range("K969").Value = 0
$ws.Range("L969").Value = 38
$ws.Range("M969").Value = 0
$ws.Range("O969").Value = 0
$ws.Range("P969").Value = 1
$ws.Range("Q969").Value = 1
$ws.Range("T969").Value = 139

# Row 970
$ws.Range("A970").Value = "'" + $dates[6]
$ws.Range("A970").Style = "Normal"
$ws.Range("B970").Value = 969
$ws.Range("C970").Value = 39486
$ws.Range("D970").Value = 691
$ws.Range("F970").Value = 0.01749987337
$ws.Range("G970").Value = 38709
$ws.Range("H970").Value = 27040
$ws.Range("I970").Value = 66526
$ws.Range("J970").Value = 38
$ws.Range("K970").Value = 0
$ws.Range("L970").Value = 38
$ws.Range("M970").Value = 0
$ws.Range("O970").Value = 0
$ws.Range("P970").Value = 1
$ws.Range("Q970").Value = 1
$ws.Range("T970").Value = 139

# Row 971
$ws.Range("A971").Value = "'" + $dates[7]
$ws.Range("A971").Style = "Normal"
$ws.Range("B971").Value = 970
$ws.Range("C971").Value = 39486
$ws.Range("D971").Value = 691
$ws.Range("F971").Value = 0.01749987337
$ws.Range("G971").Value = 38709
$ws.Range("H971").Value = 27040
$ws.Range("I971").Value = 66526
$ws.Range("J971").Value = 38
$ws.Range("K971").Value = 0
$ws.Range("L971").Value = 38
$ws.Range("M971").Value = 0
$ws.Range("O971").Value = 0
$ws.Range("P971").Value = 1
$ws.Range("Q971").Value = 1
$ws.Range("T971").Value = 139

# Row 972
$ws.Range("A972").Value = "'" + $dates[8]
$ws.Range("A972").Style = "Normal"
$ws.Range("B972").Value = 971
$ws.Range("C972").Value = 39486
$ws.Range("D972").Value = 691
$ws.Range("F972").Value = 0.01749987337
$ws.Range("G972").Value = 38709
$ws.Range("H972").Value = 27040
$ws.Range("I972").Value = 66526
$ws.Range("J972").Value = 38
$ws.Range("K972").Value = 0
$ws.Range("L972").Value = 38
$ws.Range("M972").Value = 0
$ws.Range("O972").Value = 0
$ws.Range("P972").Value = 1
$ws.Range("Q972").Value = 1
$ws.Range("T972").Value = 139

# Row 973
$ws.Range("A973").Value = "'" + $dates[9]
$ws.Range("A973").Style = "Normal"
$ws.Range("B973").Value = 972
$ws.Range("C973").Value = 39486
$ws.Range("D973").Value = 691
$ws.Range("F973").Value = 0.01749987337
$ws.Range("G973").Value = 38709
$ws.Range("H973").Value = 27040
$ws.Range("I973").Value = 66526
$ws.Range("J973").Value = 38
$ws.Range("K973").Value = 0
$ws.Range("L973").Value = 38
$ws.Range("M973").Value = 0
$ws.Range("O973").Value = 0
$ws.Range("P973").Value = 1
$ws.Range("Q973").Value = 1
$ws.Range("T973").Value = 139

# Row 974
$ws.Range("A974").Value = "'" + $dates[10]
$ws.Range("A974").Style = "Normal"
$ws.Range("B974").Value = 973
$ws.Range("C974").Value = 39486
$ws.Range("D974").Value = 691
$ws.Range("F974").Value = 0.01749987337
$ws.Range("G974").Value = 38709
$ws.Range("H974").Value = 27040
$ws.Range("I974").Value = 66526
$ws.Range("J974").Value = 38
$ws.Range("K974").Value = 0
$ws.Range("L974").Value = 38
$ws.Range("M974").Value = 0
$ws.Range("O974").Value = 0
$ws.Range("P974").Value = 1
$ws.Range("Q974").Value = 1
$ws.Range("T974").Value = 140

# Row 975
$ws.Range("A975").Value = "'" + $dates[11]
$ws.Range("A975").Style = "Normal"
$ws.Range("B975").Value = 974
$ws.Range("C975").Value = 39486
$ws.Range("D975").Value = 691
$ws.Range("F975").Value = 0.01749987337
$ws.Range("G975").Value = 38709
$ws.Range("H975").Value = 27040
$ws.Range("I975").Value = 66526
$ws.Range("J975").Value = 38
$ws.Range("K975").Value = 0
$ws.Range("L975").Value = 38
$ws.Range("M975").Value = 0
$ws.Range("O975").Value = 0
$ws.Range("P975").Value = 1
$ws.Range("Q975").Value = 1
$ws.Range("T975").Value = 140

# Row 976
$ws.Range("A976").Value = "'" + $dates[12]
$ws.Range("A976").Style = "Normal"
$ws.Range("B976").Value = 975
$ws.Range("C976").Value = 39486
$ws.Range("D976").Value = 691
$ws.Range("F976").Value = 0.01749987337
$ws.Range("G976").Value = 38709
$ws.Range("H976").Value = 27040
$ws.Range("I976").Value = 66526
$ws.Range("J976").Value = 38
$ws.Range("K976").Value = 0
$ws.Range("L976").Value = 38
$ws.Range("M976").Value = 0
$ws.Range("O976").Value = 0
$ws.Range("P976").Value = 1
$ws.Range("Q976").Value = 1
$ws.Range("T976").Value = 140

# Row 977
$ws.Range("A977").Value = "'" + $dates[13]
$ws.Range("A977").Style = "Normal"
$ws.Range("B977").Value = 976
$ws.Range("C977").Value = 39486
$ws.Range("D977").Value = 691
$ws.Range("F977").Value = 0.01749987337
$ws.Range("G977").Value = 38709
$ws.Range("H977").Value = 27040
$ws.Range("I977").Value = 66526
$ws.Range("J977").Value = 38
$ws.Range("K977").Value = 0
$ws.Range("L977").Value = 38
$ws.Range("M977").Value = 0
$ws.Range("O977").Value = 0
$ws.Range("P977").Value = 1
$ws.Range("Q977").Value = 1
$ws.Range("T977").Value = 140

# Row 978
$ws.Range("A978").Value = "'" + $dates[14]
$ws.Range("A978").Style = "Normal"
$ws.Range("B978").Value = 977
$ws.Range("C978").Value = 39882
$ws.Range("D978").Value = 693
$ws.Range("F978").Value = 0.01737625997
$ws.Range("G978").Value = 38980
$ws.Range("H978").Value = 27040
$ws.Range("I978").Value = 66922
$ws.Range("J978").Value = 396
$ws.Range("K978").Value = 2
$ws.Range("L978").Value = 396
$ws.Range("M978").Value = 1
$ws.Range("O978").Value = 1
$ws.Range("P978").Value = 5
$ws.Range("Q978").Value = 0
$ws.Range("T978").Value = 140

# Row 979
$ws.Range("A979").Value = "'" + $dates[15]
$ws.Range("A979").Style = "Normal"
$ws.Range("B979").Value = 978
$ws.Range("C979").Value = 39882
$ws.Range("D979").Value = 693
$ws.Range("F979").Value = 0.01737625997
$ws.Range("G979").Value = 38980
$ws.Range("H979").Value = 27040
$ws.Range("I979").Value = 66922
$ws.Range("J979").Value = 396
$ws.Range("K979").Value = 2
$ws.Range("L979").Value = 396
$ws.Range("M979").Value = 1
$ws.Range("O979").Value = 1
$ws.Range("P979").Value = 5
$ws.Range("Q979").Value = 0
$ws.Range("T979").Value = 140

# Row 980
$ws.Range("A980").Value = "'" + $dates[16]
$ws.Range("A980").Style = "Normal"
$ws.Range("B980").Value = 979
$ws.Range("C980").Value = 39882
$ws.Range("D980").Value = 693
$ws.Range("F980").Value = 0.01737625997
$ws.Range("G980").Value = 38980
$ws.Range("H980").Value = 27040
$ws.Range("I980").Value = 66922
$ws.Range("J980").Value = 396
$ws.Range("K980").Value = 2
$ws.Range("L980").Value = 396
$ws.Range("M980").Value = 1
$ws.Range("O980").Value = 1
$ws.Range("P980").Value = 5
$ws.Range("Q980").Value = 0
$ws.Range("T980").Value = 140

# Row 981
$ws.Range("A981").Value = "'" + $dates[17]
$ws.Range("A981").Style = "Normal"
$ws.Range("B981").Value = 980
$ws.Range("C981").Value = 39882
$ws.Range("D981").Value = 693
$ws.Range("F981").Value = 0.01737625997
$ws.Range("G981").Value = 38980
$ws.Range("H981").Value = 27040
$ws.Range("I981").Value = 66922
$ws.Range("J981").Value = 396
$ws.Range("K981").Value = 2
$ws.Range("L981").Value = 396
$ws.Range("M981").Value = 1
$ws.Range("O981").Value = 1
$ws.Range("P981").Value = 5
$ws.Range("Q981").Value = 0
$ws.Range("T981").Value = 141

# Row 982
$ws.Range("A982").Value = "'" + $dates[18]
$ws.Range("A982").Style = "Normal"
$ws.Range("B982").Value = 981
$ws.Range("C982").Value = 39882
$ws.Range("D982").Value = 693
$ws.Range("F982").Value = 0.01737625997
$ws.Range("G982").Value = 38980
$ws.Range("H982").Value = 27040
$ws.Range("I982").Value = 66922
$ws.Range("J982").Value = 396
$ws.Range("K982").Value = 2
$ws.Range("L982").Value = 396
$ws.Range("M982").Value = 1
$ws.Range("O982").Value = 1
$ws.Range("P982").Value = 5
$ws.Range("Q982").Value = 0
$ws.Range("T982").Value = 141

# Row 983
$ws.Range("A983").Value = "'" + $dates[19]
$ws.Range("A983").Style = "Normal"
$ws.Range("B983").Value = 982
$ws.Range("C983").Value = 39882
$ws.Range("D983").Value = 693
$ws.Range("F983").Value = 0.01737625997
$ws.Range("G983").Value = 38980
$ws.Range("H983").Value = 27040
$ws.Range("I983").Value = 66922
$ws.Range("J983").Value = 396
$ws.Range("K983").Value = 2
$ws.Range("L983").Value = 396
$ws.Range("M983").Value = 1
$ws.Range("O983").Value = 1
$ws.Range("P983").Value = 5
$ws.Range("Q983").Value = 0
$ws.Range("T983").Value = 141

# Row 984
$ws.Range("A984").Value = "'" + $dates[20]
$ws.Range("A984").Style = "Normal"
$ws.Range("B984").Value = 983
$ws.Range("C984").Value = 39882
$ws.Range("D984").Value = 693
$ws.Range("F984").Value = 0.01737625997
$ws.Range("G984").Value = 38980
$ws.Range("H984").Value = 27040
$ws.Range("I984").Value = 66922
$ws.Range("J984").Value = 396
$ws.Range("K984").Value = 2
$ws.Range("L984").Value = 396
$ws.Range("M984").Value = 1
$ws.Range("O984").Value = 1
$ws.Range("P984").Value = 5
$ws.Range("Q984").Value = 0
$ws.Range("T984").Value = 141

# Row 985
$ws.Range("A985").Value = "'" + $dates[21]
$ws.Range("A985").Style = "Normal"
$ws.Range("B985").Value = 984
$ws.Range("C985").Value = 39882
$ws.Range("D985").Value = 693
$ws.Range("F985").Value = 0.01737625997
$ws.Range("G985").Value = 38980
$ws.Range("H985").Value = 27040
$ws.Range("I985").Value = 66922
$ws.Range("J985").Value = 396
$ws.Range("K985").Value = 2
$ws.Range("L985").Value = 396
$ws.Range("M985").Value = 1
$ws.Range("O985").Value = 1
$ws.Range("P985").Value = 5
$ws.Range("Q985").Value = 0
$ws.Range("T985").Value = 141

# Row 986
$ws.Range("A986").Value = "'" + $dates[22]
$ws.Range("A986").Style = "Normal"
$ws.Range("B986").Value = 985
$ws.Range("C986").Value = 39882
$ws.Range("D986").Value = 693
$ws.Range("F986").Value = 0.01737625997
$ws.Range("G986").Value = 38980
$ws.Range("H986").Value = 27040
$ws.Range("I986").Value = 66922
$ws.Range("J986").Value = 396
$ws.Range("K986").Value = 2
$ws.Range("L986").Value = 396
$ws.Range("M986").Value = 1
$ws.Range("O986").Value = 1
$ws.Range("P986").Value = 5
$ws.Range("Q986").Value = 0
$ws.Range("T986").Value = 141

# Row 987
$ws.Range("A987").Value = "'" + $dates[23]
$ws.Range("A987").Style = "Normal"
$ws.Range("B987").Value = 986
$ws.Range("C987").Value = 39882
$ws.Range("D987").Value = 693
$ws.Range("F987").Value = 0.01737625997
$ws.Range("G987").Value = 38980
$ws.Range("H987").Value = 27040
$ws.Range("I987").Value = 66922
$ws.Range("J987").Value = 396
$ws.Range("K987").Value = 2
$ws.Range("L987").Value = 396
$ws.Range("M987").Value = 1
$ws.Range("O987").Value = 1
$ws.Range("P987").Value = 5
$ws.Range("Q987").Value = 0
$ws.Range("T987").Value = 141

# Row 988
$ws.Range("A988").Value = "'" + $dates[24]
$ws.Range("A988").Style = "Normal"
$ws.Range("B988").Value = 987
$ws.Range("C988").Value = 39882
$ws.Range("D988").Value = 693
$ws.Range("F988").Value = 0.01737625997
$ws.Range("G988").Value = 38980
$ws.Range("H988").Value = 27040
$ws.Range("I988").Value = 66922
$ws.Range("J988").Value = 396
$ws.Range("K988").Value = 2
$ws.Range("L988").Value = 396
$ws.Range("M988").Value = 1
$ws.Range("O988").Value = 1
$ws.Range("P988").Value = 5
$ws.Range("Q988").Value = 0
$ws.Range("T988").Value = 142

# Row 989
$ws.Range("A989").Value = "'" + $dates[25]
$ws.Range("A989").Style = "Normal"
$ws.Range("B989").Value = 988
$ws.Range("C989").Value = 39882
$ws.Range("D989").Value = 693
$ws.Range("F989").Value = 0.01737625997
$ws.Range("G989").Value = 38980
$ws.Range("H989").Value = 27040
$ws.Range("I989").Value = 66922
$ws.Range("J989").Value = 396
$ws.Range("K989").Value = 2
$ws.Range("L989").Value = 396
$ws.Range("M989").Value = 1
$ws.Range("O989").Value = 1
$ws.Range("P989").Value = 5
$ws.Range("Q989").Value = 0
$ws.Range("T989").Value = 142

# Row 990
$ws.Range("A990").Value = "'" + $dates[26]
$ws.Range("A990").Style = "Normal"
$ws.Range("B990").Value = 989
$ws.Range("C990").Value = 39882
$ws.Range("D990").Value = 693
$ws.Range("F990").Value = 0.01737625997
$ws.Range("G990").Value = 38980
$ws.Range("H990").Value = 27040
$ws.Range("I990").Value = 66922
$ws.Range("J990").Value = 396
$ws.Range("K990").Value = 2
$ws.Range("L990").Value = 396
$ws.Range("M990").Value = 1
$ws.Range("O990").Value = 1
$ws.Range("P990").Value = 5
$ws.Range("Q990").Value = 0
$ws.Range("T990").Value = 142

# Row 991
$ws.Range("A991").Value = "'" + $dates[27]
$ws.Range("A991").Style = "Normal"
$ws.Range("B991").Value = 990
$ws.Range("C991").Value = 39882
$ws.Range("D991").Value = 693
$ws.Range("F991").Value = 0.01737625997
$ws.Range("G991").Value = 38980
$ws.Range("H991").Value = 27040
$ws.Range("I991").Value = 66922
$ws.Range("J991").Value = 396
$ws.Range("K991").Value = 2
$ws.Range("L991").Value = 396
$ws.Range("M991").Value = 1
$ws.Range("O991").Value = 1
$ws.Range("P991").Value = 5
$ws.Range("Q991").Value = 0
$ws.Range("T991").Value = 142

# Row 992
$ws.Range("A992").Value = "'" + $dates[28]
$ws.Range("A992").Style = "Normal"
$ws.Range("B992").Value = 991
$ws.Range("C992").Value = 39882
$ws.Range("D992").Value = 693
$ws.Range("F992").Value = 0.01737625997
$ws.Range("G992").Value = 38980
$ws.Range("H992").Value = 27040
$ws.Range("I992").Value = 66922
$ws.Range("J992").Value = 396
$ws.Range("K992").Value = 2
$ws.Range("L992").Value = 396
$ws.Range("M992").Value = 1
$ws.Range("O992").Value = 1
$ws.Range("P992").Value = 5
$ws.Range("Q992").Value = 0
$ws.Range("T992").Value = 142

# Row 993
$ws.Range("A993").Value = "'" + $dates[29]
$ws.Range("A993").Style = "Normal"
$ws.Range("B993").Value = 992
$ws.Range("C993").Value = 40509
$ws.Range("D993").Value = 698
$ws.Range("F993").Value = 0.01723073885
$ws.Range("G993").Value = 39660
$ws.Range("H993").Value = 27040
$ws.Range("I993").Value = 67549
$ws.Range("J993").Value = 627
$ws.Range("K993").Value = 5
$ws.Range("L993").Value = 627
$ws.Range("M993").Value = 1
$ws.Range("O993").Value = 1
$ws.Range("P993").Value = 3
$ws.Range("Q993").Value = 3
$ws.Range("T993").Value = 142
